# Update cached formula results / server results in the 6 year-sheets
# (2025, 2030, 2035, 2040, 2045, 2050) of the workbook as produced by the
# latest run on the server.

$wb = $excel.ActiveWorkbook

$ws2025 = $wb.Worksheets.Item("2025")
$ws2025.Range("N2").Value = 7155.075790473336
$ws2025.Range("O2").Value = 6980.325566461754

$ws2030 = $wb.Worksheets.Item("2030")
$ws2030.Range("B2").Value = 5707.815717280662
$ws2030.Range("I2").Value = 44492.05901988943
$ws2030.Range("L2").Value = 66334.06707325629
$ws2030.Range("M2").Value = 21991.42050229464
$ws2030.Range("O2").Value = 12079.40905079305

$ws2035 = $wb.Worksheets.Item("2035")
$ws2035.Range("A2").Value = 2927.360317916481
$ws2035.Range("B2").Value = 7940.887964949257
$ws2035.Range("E2").Value = 67179.99183625776
$ws2035.Range("I2").Value = 59530.75343380851
$ws2035.Range("L2").Value = 66334.06707325629
$ws2035.Range("M2").Value = 25547.11936466757
$ws2035.Range("N2").Value = 15117.91059331085
$ws2035.Range("O2").Value = 14761.05415301146

$ws2040 = $wb.Worksheets.Item("2040")
$ws2040.Range("A2").Value = 2927.360317916481
$ws2040.Range("B2").Value = 7940.887964949257
$ws2040.Range("E2").Value = 67179.99183625776
$ws2040.Range("I2").Value = 59530.75343380851
$ws2040.Range("L2").Value = 66334.06707325629
$ws2040.Range("M2").Value = 25547.11936466757
$ws2040.Range("N2").Value = 15225.0345013318
$ws2040.Range("O2").Value = 14761.05415301146

$ws2045 = $wb.Worksheets.Item("2045")
$ws2045.Range("A2").Value = 6352.985609279765
$ws2045.Range("B2").Value = 7940.887964949257
$ws2045.Range("E2").Value = 67179.99183625776
$ws2045.Range("I2").Value = 59530.75343380851
$ws2045.Range("L2").Value = 66334.06707325629
$ws2045.Range("M2").Value = 25547.11936466757
$ws2045.Range("N2").Value = 15769.76205278203
$ws2045.Range("O2").Value = 17096.52013936021

$ws2050 = $wb.Worksheets.Item("2050")
$ws2050.Range("A2").Value = 6352.985609279765
$ws2050.Range("B2").Value = 7940.887964949257
$ws2050.Range("E2").Value = 67179.99183625776
$ws2050.Range("I2").Value = 59530.75343380851
$ws2050.Range("L2").Value = 66334.06707325629
$ws2050.Range("M2").Value = 25547.11936466757
$ws2050.Range("N2").Value = 15769.76205278203
$ws2050.Range("O2").Value = 17096.52013936021
